$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.607.75"
$ws.Range("E2").Value = "  +1.71%  "
$ws.Range("D3").Value = "2.528.58"
$ws.Range("E3").Value = "  -1.81%  "
$ws.Range("D4").Value = "'" + "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'" + "591.13"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.13%  "
$ws.Range("D6").Value = "'" + "175.15"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.50%  "
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("E8").Value = "  +0.63%  "
$ws.Range("D9").Value = "2.527.69"
$ws.Range("E9").Value = "  -1.80%  "
$ws.Range("E10").Value = "  +0.65%  "
$ws.Range("E11").Value = "  +2.23%  "
$ws.Range("E12").Value = "  +0.09%  "
$ws.Range("E13").Value = "  -3.00%  "
$ws.Range("D14").Value = "'" + "26.70"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.14%  "
$ws.Range("D15").Value = "2.992.78"
$ws.Range("E15").Value = "  -1.54%  "
$ws.Range("E16").Value = "  +0.10%  "
$ws.Range("D17").Value = "67.355.31"
$ws.Range("E17").Value = "  +1.62%  "
$ws.Range("D18").Value = "2.513.46"
$ws.Range("E18").Value = "  -1.53%  "
$ws.Range("D19").Value = "'" + "8.17"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +5.37%  "
$ws.Range("D20").Value = "'" + "11.41"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.01%  "
$ws.Range("D21").Value = "'" + "357.28"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.85%  "
$ws.Range("E22").Value = "  -0.96%  "
$ws.Range("D23").Value = "'" + "4.64"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.63%  "
$ws.Range("D24").Value = "'" + "1.99"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +4.90%  "
$ws.Range("E25").Value = "  -0.01%  "
$ws.Range("D26").Value = "'" + "10.25"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +4.14%  "
$ws.Range("D27").Value = "'" + "69.89"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.10%  "
$ws.Range("D28").Value = "'" + "1.00"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.09%  "
$ws.Range("D29").Value = "2.652.18"
$ws.Range("E29").Value = "  -2.05%  "
$ws.Range("D30").Value = "'" + "0.0" + [char]8323 + "0987"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.06%  "
$ws.Range("D31").Value = "'" + "551.99"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +4.49%  "
$ws.Range("D32").Value = "'" + "8.26"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.44%  "
$ws.Range("E33").Value = "  +1.67%  "
$ws.Range("E34").Value = "  +0.67%  "
$ws.Range("D35").Value = "'" + "0.130"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.08%  "
$ws.Range("E36").Value = "  +0.00%  "
$ws.Range("E37").Value = "  +1.11%  "
$ws.Range("D38").Value = "'" + "158.18"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.38%  "
$ws.Range("D39").Value = "'" + "18.76"
$ws.Range("D39").Style = "Normal"
$ws.Range("E40").Value = "  +0.82%  "
$ws.Range("D41").Value = "'" + "0.356"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.91%  "
$ws.Range("E42").Value = "  +1.57%  "
$ws.Range("D43").Value = "'" + "5.15"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.80%  "
$ws.Range("D44").Value = "'" + "2.58"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +5.76%  "
$ws.Range("E45").Value = "  -0.01%  "
$ws.Range("D46").Value = "'" + "149.27"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.09%  "
$ws.Range("D47").Value = "'" + "0.559"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.28%  "
$ws.Range("B48").Value = "Filecoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D48").Value = "'" + "3.70"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.46%  "
$ws.Range("B49").Value = "BabyDogeCoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D49").Value = "'" + "0.0" + [char]8326 + "0276"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.36%  "
$ws.Range("E50").Value = "  -1.42%  "
$ws.Range("E51").Value = "  -0.18%  "
